$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the overdue-period label shown for every worker row: "2507" -> "2508"
#    (shared string used by E16:E22 - updating every occurrence so the shared
#    string pool collapses/re-dedupes exactly like a normal user edit would)
foreach ($r in 16..22) {
    $ws.Range("E$r").Value = "2508"
}

# 2) Update "VALOR MORA" total figure
$ws.Range("E11").Value = 170820

# 3) Update "Cant. Trabajadores" (worker count) figure
$ws.Range("C13").Value = 3

# 4) Before removing the middle worker rows, copy the closing/bottom border
#    formatting (currently on the last row of the table, row 22) onto row 18,
#    which will become the new last row of the table once rows 19:22 are
#    removed.
$ws.Range("B22:J22").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 5) Remove the 4 worker rows that are no longer part of this statement
#    (PEDRO ANTONIO MORENO GOMEZ, JUAN DAVID PINZON QUIROGA,
#     NANCY DEL SOCORRO GOMEZ TORRES, JOSE JULIAN ZUÑIGA MORENO)
$ws.Rows("19:22").Delete()

# 6) Column D ("Nombre Trabajador") no longer needs to fit the longer names
#    that were removed - shrink it back down to fit what's left.
$ws.Columns("D").ColumnWidth = 29
